$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update overall P/L percentage shown in H2
$ws.Range("H2").Value = "P/L  12.52%"

# Update BTC-USD row (row 7) figures
$ws.Range("C7").Value = 35256.523

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-9576.48"

$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = "-287.27"

$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "-21.36 %"
